# Enemy Dark Elf Design
# Adds 7 new "Dark Elf" enemy rows (rows 9-15) to the NATURE sheet (sheet3),
# matching the FIRE (Orc) sheet's layout/formula pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NATURE")

# Column layout (see header row 1): A id, B name, C element, D hp, E dmg,
# F exp, G money, H speed, I imgsrc (formula), J pivot, K pixelperunit

$rows = @(
    @{ Row=9;  A=2020100; B="다크엘프 검사, Dark Elf Soldier";        J="0.5,0.0625" },
    @{ Row=10; A=2020101; B="다크엘프 창술사, Dark Elf Spearman";      J="0.5,0.0625" },
    @{ Row=11; A=2020102; B="다크엘프 어쌔신, Dark Elf Assassin";      J="0.5,0.0625" },
    @{ Row=12; A=2021102; B="다크엘프 어쌔신, Dark Elf Assassin";      J="0.5,0.0625" },
    @{ Row=13; A=2021103; B="다크엘프 마법사, Dark Elf Mage";          J="0.5,0.0625" },
    @{ Row=14; A=2021104; B="다크엘프 마검사, Dark Elf Magic Knight";  J="0.5,0.0625" },
    @{ Row=15; A=2022105; B="다크엘프 로드, Dark Elf Lord";            J="0.5,0.0625" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = 2
    $ws.Range("D$n").Value = 2
    $ws.Range("E$n").Value = 1
    $ws.Range("F$n").Value = 2
    $ws.Range("G$n").Value = 50
    $ws.Range("H$n").Value = 5
    $ws.Range("I$n").Formula = "=CONCATENATE(""/Sprites/Enemy/"",A$n,""/"")"
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = 24
}

# Move the active selection to where Excel would land after entering the
# last row of data (one row below, first empty column past J).
[void]$ws.Range("J16").Select()
